$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$data = @{
    2  = @(0.3624364103015721, 2.902560095337955, 14.15319625081184, 3.762073397850159, 3.78183458407233, 51)
    3  = @(0.03493740399305523, 3.008365163207992, 12.2517222994939, 3.500246034137301, 3.535606300991146, 50)
    4  = @(0.376526435448383, 2.512159032405139, 10.14738472548122, 3.185495993637603, 3.195944747761078, 49)
    5  = @(0.128062596573612, 2.051014053753879, 8.415776543719975, 2.900995784850432, 2.92883708213644, 48)
    6  = @(0.3895087488094023, 2.65705602853736, 11.08899265343494, 3.330013911898108, 3.342909279152719, 47)
    7  = @(0.1522837063575818, 2.670556882397438, 10.36728199909625, 3.219826392695148, 3.251762656728451, 46)
    8  = @(0.3987337555141685, 2.763470879128723, 12.14087511611534, 3.48437585746936, 3.500600250070111, 45)
    9  = @(0.1063170287658169, 2.290045569855232, 9.722226614463244, 3.118048526637013, 3.152262477892468, 44)
    10 = @(0.458157375663734, 2.282092575471136, 8.777379630194348, 2.962664279022237, 2.961664880566207, 43)
    11 = @(0.06186012078826768, 2.779572087388781, 10.64187340981079, 3.262188438734156, 3.301137878033202, 42)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
